$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 26.144619
$ws.Range("H2").Value = 78.433857
$ws.Range("I2").Value = 0.5211737020083955
$ws.Range("J2").Value = 0.5211737020083955
$ws.Range("M2").Value = 2.724001666666667
$ws.Range("N2").Value = 8.172005
$ws.Range("O2").Value = 0.04635500474236593
$ws.Range("P2").Value = 0.04635500474236593
$ws.Range("Q2").Value = 71.217985730365
$ws.Range("R2").Value = 640.9618715732851
$ws.Range("S2").Value = 0.02415900942819558
$ws.Range("T2").Value = 0.02415900942819558

# Row 3
$ws.Range("G3").Value = 26.144619
$ws.Range("H3").Value = 78.433857
$ws.Range("I3").Value = 0.5211737020083955
$ws.Range("J3").Value = 0.5211737020083955
$ws.Range("O3").Value = 0.6912512390256352
$ws.Range("P3").Value = 0.6912512390256351
$ws.Range("Q3").Value = 1062.010912319716
$ws.Range("R3").Value = 9558.098210877446
$ws.Range("S3").Value = 0.3602619672608806
$ws.Range("T3").Value = 0.3602619672608806

# Row 4
$ws.Range("G4").Value = 26.144619
$ws.Range("H4").Value = 78.433857
$ws.Range("I4").Value = 0.5211737020083955
$ws.Range("J4").Value = 0.5211737020083955
$ws.Range("M4").Value = 15.419285
$ws.Range("N4").Value = 46.257855
$ws.Range("O4").Value = 0.2623937562319988
$ws.Range("P4").Value = 0.2623937562319988
$ws.Range("Q4").Value = 403.131331577415
$ws.Range("R4").Value = 3628.181984196735
$ws.Range("S4").Value = 0.1367527253193193
$ws.Range("T4").Value = 0.1367527253193193

# Row 5
$ws.Range("I5").Value = 0.3571392594830743
$ws.Range("J5").Value = 0.3571392594830742
$ws.Range("M5").Value = 2.724001666666667
$ws.Range("N5").Value = 8.172005
$ws.Range("O5").Value = 0.04635500474236593
$ws.Range("P5").Value = 0.04635500474236593
$ws.Range("Q5").Value = 48.80280525975001
$ws.Range("R5").Value = 439.22524733775
$ws.Range("S5").Value = 0.01655519206702297
$ws.Range("T5").Value = 0.01655519206702296

# Row 6
$ws.Range("I6").Value = 0.3571392594830743
$ws.Range("J6").Value = 0.3571392594830742
$ws.Range("O6").Value = 0.6912512390256352
$ws.Range("P6").Value = 0.6912512390256351
$ws.Range("Q6").Value = 727.7531259294002
$ws.Range("R6").Value = 6549.778133364601
$ws.Range("S6").Value = 0.246872955622373
$ws.Range("T6").Value = 0.2468729556223729

# Row 7
$ws.Range("I7").Value = 0.3571392594830743
$ws.Range("J7").Value = 0.3571392594830742
$ws.Range("M7").Value = 15.419285
$ws.Range("N7").Value = 46.257855
$ws.Range("O7").Value = 0.2623937562319988
$ws.Range("P7").Value = 0.2623937562319988
$ws.Range("Q7").Value = 276.24959716725
$ws.Range("R7").Value = 2486.24637450525
$ws.Range("S7").Value = 0.09371111179367837
$ws.Range("T7").Value = 0.09371111179367836

# Row 8
$ws.Range("G8").Value = 6.104416333333333
$ws.Range("H8").Value = 18.313249
$ws.Range("I8").Value = 0.1216870385085301
$ws.Range("J8").Value = 0.1216870385085301
$ws.Range("M8").Value = 2.724001666666667
$ws.Range("N8").Value = 8.172005
$ws.Range("O8").Value = 0.04635500474236593
$ws.Range("P8").Value = 0.04635500474236593
$ws.Range("Q8").Value = 16.62844026602722
$ws.Range("R8").Value = 149.655962394245
$ws.Range("S8").Value = 0.00564080324714738
$ws.Range("T8").Value = 0.00564080324714738

# Row 9
$ws.Range("G9").Value = 6.104416333333333
$ws.Range("H9").Value = 18.313249
$ws.Range("I9").Value = 0.1216870385085301
$ws.Range("J9").Value = 0.1216870385085301
$ws.Range("O9").Value = 0.6912512390256352
$ws.Range("P9").Value = 0.6912512390256351
$ws.Range("Q9").Value = 247.9652413119009
$ws.Range("R9").Value = 2231.687171807108
$ws.Range("S9").Value = 0.08411631614238164
$ws.Range("T9").Value = 0.08411631614238163

# Row 10
$ws.Range("G10").Value = 6.104416333333333
$ws.Range("H10").Value = 18.313249
$ws.Range("I10").Value = 0.1216870385085301
$ws.Range("J10").Value = 0.1216870385085301
$ws.Range("M10").Value = 15.419285
$ws.Range("N10").Value = 46.257855
$ws.Range("O10").Value = 0.2623937562319988
$ws.Range("P10").Value = 0.2623937562319988
$ws.Range("Q10").Value = 94.12573520232166
$ws.Range("R10").Value = 847.1316168208949
$ws.Range("S10").Value = 0.03192991911900112
$ws.Range("T10").Value = 0.03192991911900112
